$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.701.42'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.756.00'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.00'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.754.01'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.29'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.379.34'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.751.94'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.727.10'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.29'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.58%  '
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.88'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +19.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '495.80'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.731'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000151'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.99%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.33'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.22'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.53'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +7.18%  '
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.95'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.08'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.900.45'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.689.80'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.02'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '445.39'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.02'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('E45').Value = '  +3.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.49'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.43%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.36'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.845.77'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.84'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.14%  '
$ws.Range('E51').Value = '  +2.86%  '
